# Generate Report for Handback
# ------------------------------------------------------------------
# This script applies the "handback" updates to the localization
# status report:
#   * Overview / zh-cn / de-de "Status" cells flip from
#     "Ready for handoff" to "Handed back: in sync with en-US"
#   * zh-cn / de-de rows gain a "Latest Target File" hyperlink,
#     a "Latest Handback File" name and an updated
#     "Latest Handback DateTime" timestamp
#   * Columns that now hold longer text are widened to fit
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$mdFileName   = "de2f0e6c-f087-4fa4-9c6f-c804fd303f52.md"
$mdHyperlink  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e7702cb3c6500e373bfebf1d07f40f4a60a6f47d/e2e/" + $mdFileName
$statusText   = "Handed back: in sync with en-US"

$zhXlfFileName = "de2f0e6c-f087-4fa4-9c6f-c804fd303f52.1fc4f996a9086de4e1176e20ec32355e52f818d3.zh-cn.xlf"
$deXlfFileName = "de2f0e6c-f087-4fa4-9c6f-c804fd303f52.1fc4f996a9086de4e1176e20ec32355e52f818d3.de-de.xlf"

$zhHandbackDateTime = "2016-08-18 00:55:09"
$deHandbackDateTime = "2016-08-18 00:55:17"

# ------------------------------------------------------------------
# Overview sheet: Status columns for both languages
# ------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Columns.Item(5).ColumnWidth = 29.085
$wsOverview.Columns.Item(6).ColumnWidth = 29.085

# ------------------------------------------------------------------
# zh-cn sheet
# ------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $statusText

$wsZh.Range("I2").Value = $mdFileName
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $mdHyperlink, "", "", $mdFileName) | Out-Null
$wsZh.Range("I2").Style = "HyperLink"

$wsZh.Range("J2").Value = $zhXlfFileName
$wsZh.Range("K2").Value = $zhHandbackDateTime

$wsZh.Columns.Item(3).ColumnWidth = 29.085
$wsZh.Columns.Item(9).ColumnWidth = 38.42
$wsZh.Columns.Item(10).ColumnWidth = 39.085

# ------------------------------------------------------------------
# de-de sheet
# ------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $statusText

$wsDe.Range("I2").Value = $mdFileName
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $mdHyperlink, "", "", $mdFileName) | Out-Null
$wsDe.Range("I2").Style = "HyperLink"

$wsDe.Range("J2").Value = $deXlfFileName
$wsDe.Range("K2").Value = $deHandbackDateTime

$wsDe.Columns.Item(3).ColumnWidth = 29.085
$wsDe.Columns.Item(9).ColumnWidth = 38.42
$wsDe.Columns.Item(10).ColumnWidth = 39.085
